# Realestate Update resale numbers 2025-01-14 22:42
# Appends a new data row (row 25) to the CityResaleNum sheet with the
# 2025-01-14 22:42:40 resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Columns A, B, D hold values that look numeric/date-like but must stay as
# plain text (e.g. "02" keeps its leading zero). Force a text number format
# before assigning, then strip the format back off so the cell ends up with
# no explicit style, matching the rest of the data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-14"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "22:42:40"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric resale-count columns (Beijing .. Wuhan)
$ws.Cells.Item($row, 5).Value = 126797
$ws.Cells.Item($row, 6).Value = 143477
$ws.Cells.Item($row, 7).Value = 169344
$ws.Cells.Item($row, 8).Value = 155438
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142850
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192998
$ws.Cells.Item($row, 14).Value = 115423
$ws.Cells.Item($row, 15).Value = 46011
$ws.Cells.Item($row, 16).Value = 28511
$ws.Cells.Item($row, 17).Value = 65643
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49249
$ws.Cells.Item($row, 20).Value = -1
